# "Finished Task 5 (Email Send)" - the email setup notes/links and the
# "Admin Login Area To Set Price?" idea are no longer needed on the sheet,
# idea #5 is now marked DONE, and the summary blurb at the top (F2) is
# swapped out for the description of the smart-parking system.
#
# Note: cell contents are cleared in place (not via row/column Delete) so
# that row numbers below (e.g. the styled D14 cell) do not shift - this
# matches the target layout where dimension shrinks to B2:F14 but row 14
# keeps its original row number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: remove "Admin Login Area To Set Price?"
$ws.Range("F12").ClearContents()

# Rows 16-19: remove the "Email Guide Links for 5" section
#   E16 "Email Guide Links for 5"
#   E17 "How To Setup Email Function PHP: ..."
#   E18 "How To Setup gmail SMTP server form PHP: ..."
#   E19 "I May add or remove things from this assignment ..."
$ws.Range("E16:E19").ClearContents()

# F2: swap the old note for the new project summary text
$ws.Range("F2").Value = "This smart paking system is designed for a pre setup group of people such a hotell where you have users details"

# Idea #5 (row 10) is now finished -> mark its Status column as DONE
$ws.Range("C10").Value = "DONE"

# Update the saved view: scroll/selection moved to M8
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("M8").Select()
